# Applies the dataset re-run: two new measurement rows ("Holden", "Rizzie Spiral")
# inserted after "Spiral5", all subsequent rows shift down by two, two more rows
# ("Michael-CCHex", "Michael-SNHex") appended at the end, and "Thomas Hex" was
# renamed to "Matthies Hex" (now the row for index 9, at sheet row 11).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Cells.Item(4, 2).Value = 'Holden'
$ws.Cells.Item(4, 3).Value = 0.9228545015898802
$ws.Cells.Item(4, 4).Value = 1.386215645924811
$ws.Cells.Item(4, 5).Value = 0.889316700598519
$ws.Cells.Item(4, 6).Value = 1.386215645924811
$ws.Cells.Item(4, 7).Value = 0.889316700598519
$ws.Cells.Item(4, 8).Value = 0.8393583130099534
$ws.Cells.Item(4, 9).Value = 1.176319683174596
$ws.Cells.Item(4, 10).Value = 0.9130295872900184
$ws.Cells.Item(4, 11).Value = 0.889316700598519
$ws.Cells.Item(4, 12).Value = 0.9228545015898802
$ws.Cells.Item(4, 13).Value = 1.154535073757346
$ws.Cells.Item(4, 14).Value = 1.154535073757346
$ws.Cells.Item(4, 15).Value = 1.161796610229763
$ws.Cells.Item(4, 16).Value = 1.06612894937107
$ws.Cells.Item(4, 17).Value = 1.06612894937107
$ws.Cells.Item(4, 18).Value = 1.021925887177932
$ws.Cells.Item(4, 19).Value = 1.021925887177932
$ws.Cells.Item(4, 20).Value = 1.02118240526463

# Row 5
$ws.Cells.Item(5, 2).Value = 'Rizzie Spiral'
$ws.Cells.Item(5, 3).Value = 1.130982459192274
$ws.Cells.Item(5, 4).Value = 0.8344367663395904
$ws.Cells.Item(5, 5).Value = 0.9041808864503452
$ws.Cells.Item(5, 6).Value = 0.8344367663395904
$ws.Cells.Item(5, 7).Value = 0.9041808864503452
$ws.Cells.Item(5, 8).Value = 1.467925981547674
$ws.Cells.Item(5, 9).Value = 0.8180078256798523
$ws.Cells.Item(5, 10).Value = 1.063368821033031
$ws.Cells.Item(5, 11).Value = 0.9041808864503452
$ws.Cells.Item(5, 12).Value = 1.130982459192274
$ws.Cells.Item(5, 13).Value = 0.9827096127659322
$ws.Cells.Item(5, 14).Value = 0.9827096127659322
$ws.Cells.Item(5, 15).Value = 0.9278090170705723
$ws.Cells.Item(5, 16).Value = 0.9565333706607365
$ws.Cells.Item(5, 17).Value = 0.9565333706607365
$ws.Cells.Item(5, 18).Value = 0.9434452496081387
$ws.Cells.Item(5, 19).Value = 0.9434452496081387
$ws.Cells.Item(5, 20).Value = 1.036483790040461

# Row 6
$ws.Cells.Item(6, 2).Value = 'RotRing OmegaMax-90'
$ws.Cells.Item(6, 3).Value = 0.9957856217843082
$ws.Cells.Item(6, 4).Value = 0.9017874889788284
$ws.Cells.Item(6, 5).Value = 1.085026010233661
$ws.Cells.Item(6, 6).Value = 0.9017874889788284
$ws.Cells.Item(6, 7).Value = 1.085026010233661
$ws.Cells.Item(6, 8).Value = 0.9958204495288052
$ws.Cells.Item(6, 9).Value = 0.9532010507387216
$ws.Cells.Item(6, 10).Value = 1.021543206145053
$ws.Cells.Item(6, 11).Value = 1.085026010233661
$ws.Cells.Item(6, 12).Value = 0.9957856217843082
$ws.Cells.Item(6, 13).Value = 0.9487865553815683
$ws.Cells.Item(6, 14).Value = 0.9487865553815683
$ws.Cells.Item(6, 15).Value = 0.9502580538339528
$ws.Cells.Item(6, 16).Value = 0.9941997069989327
$ws.Cells.Item(6, 17).Value = 0.9941997069989327
$ws.Cells.Item(6, 18).Value = 1.016906282807615
$ws.Cells.Item(6, 19).Value = 1.016906282807615
$ws.Cells.Item(6, 20).Value = 0.9921939712348963

# Row 7
$ws.Cells.Item(7, 2).Value = 'Equal Angle'
$ws.Cells.Item(7, 3).Value = 0.9081823531628253
$ws.Cells.Item(7, 4).Value = 1.449339625583576
$ws.Cells.Item(7, 5).Value = 0.8788915965778105
$ws.Cells.Item(7, 6).Value = 1.449339625583576
$ws.Cells.Item(7, 7).Value = 0.8788915965778105
$ws.Cells.Item(7, 8).Value = 0.8151365083638326
$ws.Cells.Item(7, 9).Value = 1.201591782118154
$ws.Cells.Item(7, 10).Value = 0.8995313213832831
$ws.Cells.Item(7, 11).Value = 0.8788915965778105
$ws.Cells.Item(7, 12).Value = 0.9081823531628253
$ws.Cells.Item(7, 13).Value = 1.178760989373201
$ws.Cells.Item(7, 14).Value = 1.178760989373201
$ws.Cells.Item(7, 15).Value = 1.186371253621519
$ws.Cells.Item(7, 16).Value = 1.07880452510807
$ws.Cells.Item(7, 17).Value = 1.07880452510807
$ws.Cells.Item(7, 18).Value = 1.028826292975505
$ws.Cells.Item(7, 19).Value = 1.028826292975505
$ws.Cells.Item(7, 20).Value = 1.025445531198247

# Row 8
$ws.Cells.Item(8, 2).Value = 'Tilt Rotate'
$ws.Cells.Item(8, 3).Value = 0.6868740243434884
$ws.Cells.Item(8, 4).Value = 2.656084637647103
$ws.Cells.Item(8, 5).Value = 0.5568101562590696
$ws.Cells.Item(8, 6).Value = 2.656084637647103
$ws.Cells.Item(8, 7).Value = 0.5568101562590696
$ws.Cells.Item(8, 8).Value = 0.4484186583378973
$ws.Cells.Item(8, 9).Value = 1.679307984522861
$ws.Cells.Item(8, 10).Value = 0.6483641325109282
$ws.Cells.Item(8, 11).Value = 0.5568101562590696
$ws.Cells.Item(8, 12).Value = 0.6868740243434884
$ws.Cells.Item(8, 13).Value = 1.671479330995296
$ws.Cells.Item(8, 14).Value = 1.671479330995296
$ws.Cells.Item(8, 15).Value = 1.674088882171151
$ws.Cells.Item(8, 16).Value = 1.299922939416554
$ws.Cells.Item(8, 17).Value = 1.299922939416554
$ws.Cells.Item(8, 18).Value = 1.114144743627183
$ws.Cells.Item(8, 19).Value = 1.114144743627183
$ws.Cells.Item(8, 20).Value = 1.112643265603558

# Row 9
$ws.Cells.Item(9, 2).Value = 'CLR'
$ws.Cells.Item(9, 3).Value = 0.9869408235375077
$ws.Cells.Item(9, 4).Value = 1.03848699472462
$ws.Cells.Item(9, 5).Value = 0.9985235759045246
$ws.Cells.Item(9, 6).Value = 1.03848699472462
$ws.Cells.Item(9, 7).Value = 0.9985235759045246
$ws.Cells.Item(9, 8).Value = 0.9706543668790255
$ws.Cells.Item(9, 9).Value = 1.020251874081768
$ws.Cells.Item(9, 10).Value = 0.9903031962845457
$ws.Cells.Item(9, 11).Value = 0.9985235759045246
$ws.Cells.Item(9, 12).Value = 0.9869408235375077
$ws.Cells.Item(9, 13).Value = 1.012713909131064
$ws.Cells.Item(9, 14).Value = 1.012713909131064
$ws.Cells.Item(9, 15).Value = 1.015226564114632
$ws.Cells.Item(9, 16).Value = 1.007983798055551
$ws.Cells.Item(9, 17).Value = 1.007983798055551
$ws.Cells.Item(9, 18).Value = 1.005618742517794
$ws.Cells.Item(9, 19).Value = 1.005618742517794
$ws.Cells.Item(9, 20).Value = 1.000860138568665

# Row 10
$ws.Cells.Item(10, 2).Value = 'Rizzie Hex'
$ws.Cells.Item(10, 3).Value = 0.9990378467909614
$ws.Cells.Item(10, 4).Value = 1.003866385975271
$ws.Cells.Item(10, 5).Value = 0.9994037099562305
$ws.Cells.Item(10, 6).Value = 1.003866385975271
$ws.Cells.Item(10, 7).Value = 0.9994037099562305
$ws.Cells.Item(10, 8).Value = 0.9981248687760244
$ws.Cells.Item(10, 9).Value = 1.001680976504458
$ws.Cells.Item(10, 10).Value = 0.9991457905094707
$ws.Cells.Item(10, 11).Value = 0.9994037099562305
$ws.Cells.Item(10, 12).Value = 0.9990378467909614
$ws.Cells.Item(10, 13).Value = 1.001452116383116
$ws.Cells.Item(10, 14).Value = 1.001452116383116
$ws.Cells.Item(10, 15).Value = 1.00152840309023
$ws.Cells.Item(10, 16).Value = 1.000769314240821
$ws.Cells.Item(10, 17).Value = 1.000769314240821
$ws.Cells.Item(10, 18).Value = 1.000427913169673
$ws.Cells.Item(10, 19).Value = 1.000427913169673
$ws.Cells.Item(10, 20).Value = 1.000209929752069

# Row 11
$ws.Cells.Item(11, 2).Value = 'Matthies Hex'
$ws.Cells.Item(11, 3).Value = 0.9750761903524402
$ws.Cells.Item(11, 4).Value = 1.072908334261707
$ws.Cells.Item(11, 5).Value = 0.9979270591049408
$ws.Cells.Item(11, 6).Value = 1.072908334261707
$ws.Cells.Item(11, 7).Value = 0.9979270591049408
$ws.Cells.Item(11, 8).Value = 0.9455176546297438
$ws.Cells.Item(11, 9).Value = 1.037798343618565
$ws.Cells.Item(11, 10).Value = 0.9817120679331135
$ws.Cells.Item(11, 11).Value = 0.9979270591049408
$ws.Cells.Item(11, 12).Value = 0.9750761903524402
$ws.Cells.Item(11, 13).Value = 1.023992262307073
$ws.Cells.Item(11, 14).Value = 1.023992262307073
$ws.Cells.Item(11, 15).Value = 1.028594289410904
$ws.Cells.Item(11, 16).Value = 1.015303861239696
$ws.Cells.Item(11, 17).Value = 1.015303861239696
$ws.Cells.Item(11, 18).Value = 1.010959660706007
$ws.Cells.Item(11, 19).Value = 1.010959660706007
$ws.Cells.Item(11, 20).Value = 1.001823274983418

# Row 12
$ws.Cells.Item(12, 2).Value = 'Tilt Rotate_Partial'
$ws.Cells.Item(12, 3).Value = 0.6794186276487614
$ws.Cells.Item(12, 4).Value = 2.693244748941595
$ws.Cells.Item(12, 5).Value = 0.5519497539580631
$ws.Cells.Item(12, 6).Value = 2.693244748941595
$ws.Cells.Item(12, 7).Value = 0.5519497539580631
$ws.Cells.Item(12, 8).Value = 0.4474831781199314
$ws.Cells.Item(12, 9).Value = 1.688620914273127
$ws.Cells.Item(12, 10).Value = 0.6416531588629021
$ws.Cells.Item(12, 11).Value = 0.5519497539580631
$ws.Cells.Item(12, 12).Value = 0.6794186276487614
$ws.Cells.Item(12, 13).Value = 1.686331688295178
$ws.Cells.Item(12, 14).Value = 1.686331688295178
$ws.Cells.Item(12, 15).Value = 1.687094763621161
$ws.Cells.Item(12, 16).Value = 1.308204376849473
$ws.Cells.Item(12, 17).Value = 1.308204376849473
$ws.Cells.Item(12, 18).Value = 1.119140721126621
$ws.Cells.Item(12, 19).Value = 1.119140721126621
$ws.Cells.Item(12, 20).Value = 1.11706173030073

# Row 13
$ws.Cells.Item(13, 2).Value = 'RotRing OmegaMax-60'
$ws.Cells.Item(13, 3).Value = 0.9551777475594478
$ws.Cells.Item(13, 4).Value = 0.9700286599781373
$ws.Cells.Item(13, 5).Value = 1.106333785863047
$ws.Cells.Item(13, 6).Value = 0.9700286599781373
$ws.Cells.Item(13, 7).Value = 1.106333785863047
$ws.Cells.Item(13, 8).Value = 0.8597543354199119
$ws.Cells.Item(13, 9).Value = 1.012703503913756
$ws.Cells.Item(13, 10).Value = 0.9988573581130877
$ws.Cells.Item(13, 11).Value = 1.106333785863047
$ws.Cells.Item(13, 12).Value = 0.9551777475594478
$ws.Cells.Item(13, 13).Value = 0.9626032037687926
$ws.Cells.Item(13, 14).Value = 0.9626032037687926
$ws.Cells.Item(13, 15).Value = 0.9793033038171138
$ws.Cells.Item(13, 16).Value = 1.010513397800211
$ws.Cells.Item(13, 17).Value = 1.010513397800211
$ws.Cells.Item(13, 18).Value = 1.03446849481592
$ws.Cells.Item(13, 19).Value = 1.03446849481592
$ws.Cells.Item(13, 20).Value = 0.983809231807898

# Row 14
$ws.Cells.Item(14, 2).Value = 'Equal Angle_Partial'
$ws.Cells.Item(14, 3).Value = 0.9054629497578945
$ws.Cells.Item(14, 4).Value = 1.504811564547371
$ws.Cells.Item(14, 5).Value = 0.8597119242842107
$ws.Cells.Item(14, 6).Value = 1.504811564547371
$ws.Cells.Item(14, 7).Value = 0.8597119242842107
$ws.Cells.Item(14, 8).Value = 0.844779967521052
$ws.Cells.Item(14, 9).Value = 1.206733396063159
$ws.Cells.Item(14, 10).Value = 0.8919659231999991
$ws.Cells.Item(14, 11).Value = 0.8597119242842107
$ws.Cells.Item(14, 12).Value = 0.9054629497578945
$ws.Cells.Item(14, 13).Value = 1.205137257152633
$ws.Cells.Item(14, 14).Value = 1.205137257152633
$ws.Cells.Item(14, 15).Value = 1.205669303456141
$ws.Cells.Item(14, 16).Value = 1.089995479529825
$ws.Cells.Item(14, 17).Value = 1.089995479529825
$ws.Cells.Item(14, 18).Value = 1.032424590718422
$ws.Cells.Item(14, 19).Value = 1.032424590718422
$ws.Cells.Item(14, 20).Value = 1.035577620895614

# Row 15
$ws.Cells.Item(15, 2).Value = 'Rizzie Hex_Partial'
$ws.Cells.Item(15, 3).Value = 1.051459653945521
$ws.Cells.Item(15, 4).Value = 0.7612610939260452
$ws.Cells.Item(15, 5).Value = 1.084925599298505
$ws.Cells.Item(15, 6).Value = 0.7612610939260452
$ws.Cells.Item(15, 7).Value = 1.084925599298505
$ws.Cells.Item(15, 8).Value = 1.168699781989065
$ws.Cells.Item(15, 9).Value = 0.857860371732878
$ws.Cells.Item(15, 10).Value = 1.061222540578913
$ws.Cells.Item(15, 11).Value = 1.084925599298505
$ws.Cells.Item(15, 12).Value = 1.051459653945521
$ws.Cells.Item(15, 13).Value = 0.9063603739357831
$ws.Cells.Item(15, 14).Value = 0.9063603739357831
$ws.Cells.Item(15, 15).Value = 0.8901937065348147
$ws.Cells.Item(15, 16).Value = 0.9658821157233571
$ws.Cells.Item(15, 17).Value = 0.9658821157233571
$ws.Cells.Item(15, 18).Value = 0.9956429866171441
$ws.Cells.Item(15, 19).Value = 0.9956429866171441
$ws.Cells.Item(15, 20).Value = 0.9975715069118212

# Row 16
$ws.Cells.Item(16, 2).Value = 'ND Single'
$ws.Cells.Item(16, 3).Value = 0.4490890699999996
$ws.Cells.Item(16, 4).Value = 3.939353100000003
$ws.Cells.Item(16, 5).Value = 0.2188589699999998
$ws.Cells.Item(16, 6).Value = 3.939353100000003
$ws.Cells.Item(16, 7).Value = 0.2188589699999998
$ws.Cells.Item(16, 8).Value = 0.05435712199999995
$ws.Cells.Item(16, 9).Value = 2.187941700000001
$ws.Cells.Item(16, 10).Value = 0.3808422
$ws.Cells.Item(16, 11).Value = 0.2188589699999998
$ws.Cells.Item(16, 12).Value = 0.4490890699999996
$ws.Cells.Item(16, 13).Value = 2.194221085000001
$ws.Cells.Item(16, 14).Value = 2.194221085000001
$ws.Cells.Item(16, 15).Value = 2.192127956666668
$ws.Cells.Item(16, 16).Value = 1.535767046666668
$ws.Cells.Item(16, 17).Value = 1.535767046666668
$ws.Cells.Item(16, 18).Value = 1.2065400275
$ws.Cells.Item(16, 19).Value = 1.2065400275
$ws.Cells.Item(16, 20).Value = 1.205073693666667

# Row 17
$ws.Cells.Item(17, 2).Value = 'RD Single'
$ws.Cells.Item(17, 3).Value = 0.8289786500000002
$ws.Cells.Item(17, 4).Value = 0.21885897
$ws.Cells.Item(17, 5).Value = 1.9813918
$ws.Cells.Item(17, 6).Value = 0.21885897
$ws.Cells.Item(17, 7).Value = 1.9813918
$ws.Cells.Item(17, 8).Value = 0.66920919
$ws.Cells.Item(17, 9).Value = 0.6454245799999999
$ws.Cells.Item(17, 10).Value = 1.1585466
$ws.Cells.Item(17, 11).Value = 1.9813918
$ws.Cells.Item(17, 12).Value = 0.8289786500000002
$ws.Cells.Item(17, 13).Value = 0.5239188100000001
$ws.Cells.Item(17, 14).Value = 0.5239188100000001
$ws.Cells.Item(17, 15).Value = 0.5644207333333333
$ws.Cells.Item(17, 16).Value = 1.00974314
$ws.Cells.Item(17, 17).Value = 1.00974314
$ws.Cells.Item(17, 18).Value = 1.252655305
$ws.Cells.Item(17, 19).Value = 1.252655305
$ws.Cells.Item(17, 20).Value = 0.9170682983333333

# Row 18
$ws.Cells.Item(18, 2).Value = 'TD Single'
$ws.Cells.Item(18, 3).Value = 0.8289786500000002
$ws.Cells.Item(18, 4).Value = 0.21885897
$ws.Cells.Item(18, 5).Value = 1.9813918
$ws.Cells.Item(18, 6).Value = 0.21885897
$ws.Cells.Item(18, 7).Value = 1.9813918
$ws.Cells.Item(18, 8).Value = 0.66920919
$ws.Cells.Item(18, 9).Value = 0.6454245799999999
$ws.Cells.Item(18, 10).Value = 1.1585466
$ws.Cells.Item(18, 11).Value = 1.9813918
$ws.Cells.Item(18, 12).Value = 0.8289786500000002
$ws.Cells.Item(18, 13).Value = 0.5239188100000001
$ws.Cells.Item(18, 14).Value = 0.5239188100000001
$ws.Cells.Item(18, 15).Value = 0.5644207333333333
$ws.Cells.Item(18, 16).Value = 1.00974314
$ws.Cells.Item(18, 17).Value = 1.00974314
$ws.Cells.Item(18, 18).Value = 1.252655305
$ws.Cells.Item(18, 19).Value = 1.252655305
$ws.Cells.Item(18, 20).Value = 0.9170682983333333

# Row 19
$ws.Cells.Item(19, 2).Value = 'Morris Single'
$ws.Cells.Item(19, 3).Value = 1.5132597
$ws.Cells.Item(19, 4).Value = 0.06559432699999999
$ws.Cells.Item(19, 5).Value = 0.7162996800000001
$ws.Cells.Item(19, 6).Value = 0.06559432699999999
$ws.Cells.Item(19, 7).Value = 0.7162996800000001
$ws.Cells.Item(19, 8).Value = 2.818143399999999
$ws.Cells.Item(19, 9).Value = 0.24865233
$ws.Cells.Item(19, 10).Value = 1.2778388
$ws.Cells.Item(19, 11).Value = 0.7162996800000001
$ws.Cells.Item(19, 12).Value = 1.5132597
$ws.Cells.Item(19, 13).Value = 0.7894270135
$ws.Cells.Item(19, 14).Value = 0.7894270135
$ws.Cells.Item(19, 15).Value = 0.6091687856666667
$ws.Cells.Item(19, 16).Value = 0.7650512356666667
$ws.Cells.Item(19, 17).Value = 0.7650512356666668
$ws.Cells.Item(19, 18).Value = 0.7528633467500001
$ws.Cells.Item(19, 19).Value = 0.7528633467500001
$ws.Cells.Item(19, 20).Value = 1.106631372833333

# Row 20
$ws.Cells.Item(20, 2).Value = 'Ring Perpendicular to ND'
$ws.Cells.Item(20, 3).Value = 0.701076531643836
$ws.Cells.Item(20, 4).Value = 1.835997022054794
$ws.Cells.Item(20, 5).Value = 0.9995425063013701
$ws.Cells.Item(20, 6).Value = 1.835997022054794
$ws.Cells.Item(20, 7).Value = 0.9995425063013701
$ws.Cells.Item(20, 8).Value = 0.337895131150685
$ws.Cells.Item(20, 9).Value = 1.441133253972603
$ws.Cells.Item(20, 10).Value = 0.7877548367123287
$ws.Cells.Item(20, 11).Value = 0.9995425063013701
$ws.Cells.Item(20, 12).Value = 0.701076531643836
$ws.Cells.Item(20, 13).Value = 1.268536776849315
$ws.Cells.Item(20, 14).Value = 1.268536776849315
$ws.Cells.Item(20, 15).Value = 1.326068935890411
$ws.Cells.Item(20, 16).Value = 1.17887202
$ws.Cells.Item(20, 17).Value = 1.17887202
$ws.Cells.Item(20, 18).Value = 1.134039641575343
$ws.Cells.Item(20, 19).Value = 1.134039641575343
$ws.Cells.Item(20, 20).Value = 1.017233213639269

# Row 21
$ws.Cells.Item(21, 2).Value = 'Ring Perpendicular to RD'
$ws.Cells.Item(21, 3).Value = 1.075012418947368
$ws.Cells.Item(21, 4).Value = 1.043442563052632
$ws.Cells.Item(21, 5).Value = 0.8787117957894734
$ws.Cells.Item(21, 6).Value = 1.043442563052632
$ws.Cells.Item(21, 7).Value = 0.8787117957894734
$ws.Cells.Item(21, 8).Value = 1.346592761
$ws.Cells.Item(21, 9).Value = 0.9135578252631579
$ws.Cells.Item(21, 10).Value = 1.016213122105263
$ws.Cells.Item(21, 11).Value = 0.8787117957894734
$ws.Cells.Item(21, 12).Value = 1.075012418947368
$ws.Cells.Item(21, 13).Value = 1.059227491
$ws.Cells.Item(21, 14).Value = 1.059227491
$ws.Cells.Item(21, 15).Value = 1.010670935754386
$ws.Cells.Item(21, 16).Value = 0.9990555925964912
$ws.Cells.Item(21, 17).Value = 0.9990555925964912
$ws.Cells.Item(21, 18).Value = 0.9689696433947368
$ws.Cells.Item(21, 19).Value = 0.9689696433947368
$ws.Cells.Item(21, 20).Value = 1.045588414359649

# Row 22
$ws.Cells.Item(22, 2).Value = 'Ring Perpendicular to TD'
$ws.Cells.Item(22, 3).Value = 1.075012418947368
$ws.Cells.Item(22, 4).Value = 1.043442563052632
$ws.Cells.Item(22, 5).Value = 0.8787117957894734
$ws.Cells.Item(22, 6).Value = 1.043442563052632
$ws.Cells.Item(22, 7).Value = 0.8787117957894734
$ws.Cells.Item(22, 8).Value = 1.346592761
$ws.Cells.Item(22, 9).Value = 0.9135578252631579
$ws.Cells.Item(22, 10).Value = 1.016213122105263
$ws.Cells.Item(22, 11).Value = 0.8787117957894734
$ws.Cells.Item(22, 12).Value = 1.075012418947368
$ws.Cells.Item(22, 13).Value = 1.059227491
$ws.Cells.Item(22, 14).Value = 1.059227491
$ws.Cells.Item(22, 15).Value = 1.010670935754386
$ws.Cells.Item(22, 16).Value = 0.9990555925964912
$ws.Cells.Item(22, 17).Value = 0.9990555925964912
$ws.Cells.Item(22, 18).Value = 0.9689696433947368
$ws.Cells.Item(22, 19).Value = 0.9689696433947368
$ws.Cells.Item(22, 20).Value = 1.045588414359649

# Row 23
$ws.Cells.Item(23, 2).Value = 'OffsetFTD'
$ws.Cells.Item(23, 3).Value = 1.165048391102379
$ws.Cells.Item(23, 4).Value = 0.6132311431620125
$ws.Cells.Item(23, 5).Value = 0.9489468280812742
$ws.Cells.Item(23, 6).Value = 0.6132311431620125
$ws.Cells.Item(23, 7).Value = 0.9489468280812742
$ws.Cells.Item(23, 8).Value = 1.391676314260136
$ws.Cells.Item(23, 9).Value = 0.7808622755914665
$ws.Cells.Item(23, 10).Value = 1.102329925782838
$ws.Cells.Item(23, 11).Value = 0.9489468280812742
$ws.Cells.Item(23, 12).Value = 1.165048391102379
$ws.Cells.Item(23, 13).Value = 0.889139767132196
$ws.Cells.Item(23, 14).Value = 0.889139767132196
$ws.Cells.Item(23, 15).Value = 0.8530472699519528
$ws.Cells.Item(23, 16).Value = 0.9090754541152221
$ws.Cells.Item(23, 17).Value = 0.9090754541152221
$ws.Cells.Item(23, 18).Value = 0.9190432976067351
$ws.Cells.Item(23, 19).Value = 0.9190432976067351
$ws.Cells.Item(23, 20).Value = 1.000349146330018

# Row 24
$ws.Cells.Item(24, 2).Value = 'OffsetATD'
$ws.Cells.Item(24, 3).Value = 0.9797230111250644
$ws.Cells.Item(24, 4).Value = 0.920854348411776
$ws.Cells.Item(24, 5).Value = 1.114500948424145
$ws.Cells.Item(24, 6).Value = 0.920854348411776
$ws.Cells.Item(24, 7).Value = 1.114500948424145
$ws.Cells.Item(24, 8).Value = 1.044247981991259
$ws.Cells.Item(24, 9).Value = 0.9377645740754055
$ws.Cells.Item(24, 10).Value = 1.018930492568054
$ws.Cells.Item(24, 11).Value = 1.114500948424145
$ws.Cells.Item(24, 12).Value = 0.9797230111250644
$ws.Cells.Item(24, 13).Value = 0.9502886797684202
$ws.Cells.Item(24, 14).Value = 0.9502886797684202
$ws.Cells.Item(24, 15).Value = 0.9461139778707487
$ws.Cells.Item(24, 16).Value = 1.005026102653662
$ws.Cells.Item(24, 17).Value = 1.005026102653662
$ws.Cells.Item(24, 18).Value = 1.032394814096283
$ws.Cells.Item(24, 19).Value = 1.032394814096283
$ws.Cells.Item(24, 20).Value = 1.002670226099284

# Row 25
$ws.Cells.Item(25, 2).Value = 'OffsetF45'
$ws.Cells.Item(25, 3).Value = 0.9294315679662116
$ws.Cells.Item(25, 4).Value = 0.6474123980439751
$ws.Cells.Item(25, 5).Value = 1.383247727733637
$ws.Cells.Item(25, 6).Value = 0.6474123980439751
$ws.Cells.Item(25, 7).Value = 1.383247727733637
$ws.Cells.Item(25, 8).Value = 0.7695134981280535
$ws.Cells.Item(25, 9).Value = 0.8916722145036236
$ws.Cells.Item(25, 10).Value = 1.062021937385172
$ws.Cells.Item(25, 11).Value = 1.383247727733637
$ws.Cells.Item(25, 12).Value = 0.9294315679662116
$ws.Cells.Item(25, 13).Value = 0.7884219830050934
$ws.Cells.Item(25, 14).Value = 0.7884219830050934
$ws.Cells.Item(25, 15).Value = 0.8228387268379368
$ws.Cells.Item(25, 16).Value = 0.9866972312479413
$ws.Cells.Item(25, 17).Value = 0.9866972312479413
$ws.Cells.Item(25, 18).Value = 1.085834855369365
$ws.Cells.Item(25, 19).Value = 1.085834855369365
$ws.Cells.Item(25, 20).Value = 0.9472165572934456

# Row 26
$ws.Cells.Item(26, 2).Value = 'OffsetA45'
$ws.Cells.Item(26, 3).Value = 1.069688241736297
$ws.Cells.Item(26, 4).Value = 1.013927747195857
$ws.Cells.Item(26, 5).Value = 0.8574694106388464
$ws.Cells.Item(26, 6).Value = 1.013927747195857
$ws.Cells.Item(26, 7).Value = 0.8574694106388464
$ws.Cells.Item(26, 8).Value = 1.165039530280905
$ws.Cells.Item(26, 9).Value = 0.9798262988859499
$ws.Cells.Item(26, 10).Value = 1.007695936101689
$ws.Cells.Item(26, 11).Value = 0.8574694106388464
$ws.Cells.Item(26, 12).Value = 1.069688241736297
$ws.Cells.Item(26, 13).Value = 1.041807994466077
$ws.Cells.Item(26, 14).Value = 1.041807994466077
$ws.Cells.Item(26, 15).Value = 1.021147429272701
$ws.Cells.Item(26, 16).Value = 0.9803617998570003
$ws.Cells.Item(26, 17).Value = 0.9803617998570003
$ws.Cells.Item(26, 18).Value = 0.9496387025524619
$ws.Cells.Item(26, 19).Value = 0.9496387025524619
$ws.Cells.Item(26, 20).Value = 1.015607860806591

# Row 27
$ws.Cells.Item(27, 2).Value = 'OffsetFRD'
$ws.Cells.Item(27, 3).Value = 1.16504839110238
$ws.Cells.Item(27, 4).Value = 0.6132311431620125
$ws.Cells.Item(27, 5).Value = 0.9489468280812741
$ws.Cells.Item(27, 6).Value = 0.6132311431620125
$ws.Cells.Item(27, 7).Value = 0.9489468280812741
$ws.Cells.Item(27, 8).Value = 1.391676314260136
$ws.Cells.Item(27, 9).Value = 0.7808622755914664
$ws.Cells.Item(27, 10).Value = 1.102329925782838
$ws.Cells.Item(27, 11).Value = 0.9489468280812741
$ws.Cells.Item(27, 12).Value = 1.16504839110238
$ws.Cells.Item(27, 13).Value = 0.8891397671321962
$ws.Cells.Item(27, 14).Value = 0.8891397671321962
$ws.Cells.Item(27, 15).Value = 0.8530472699519529
$ws.Cells.Item(27, 16).Value = 0.9090754541152221
$ws.Cells.Item(27, 17).Value = 0.9090754541152221
$ws.Cells.Item(27, 18).Value = 0.9190432976067351
$ws.Cells.Item(27, 19).Value = 0.9190432976067351
$ws.Cells.Item(27, 20).Value = 1.000349146330018

# Row 28
$ws.Cells.Item(28, 2).Value = 'OffsetARD'
$ws.Cells.Item(28, 3).Value = 0.9797230111250643
$ws.Cells.Item(28, 4).Value = 0.9208543484117762
$ws.Cells.Item(28, 5).Value = 1.114500948424145
$ws.Cells.Item(28, 6).Value = 0.9208543484117762
$ws.Cells.Item(28, 7).Value = 1.114500948424145
$ws.Cells.Item(28, 8).Value = 1.044247981991259
$ws.Cells.Item(28, 9).Value = 0.9377645740754056
$ws.Cells.Item(28, 10).Value = 1.018930492568054
$ws.Cells.Item(28, 11).Value = 1.114500948424145
$ws.Cells.Item(28, 12).Value = 0.9797230111250643
$ws.Cells.Item(28, 13).Value = 0.9502886797684202
$ws.Cells.Item(28, 14).Value = 0.9502886797684202
$ws.Cells.Item(28, 15).Value = 0.9461139778707487
$ws.Cells.Item(28, 16).Value = 1.005026102653662
$ws.Cells.Item(28, 17).Value = 1.005026102653662
$ws.Cells.Item(28, 18).Value = 1.032394814096283
$ws.Cells.Item(28, 19).Value = 1.032394814096283
$ws.Cells.Item(28, 20).Value = 1.002670226099284

# Row 29
$ws.Cells.Item(29, 2).Value = 'Gaussian Quadrature'
$ws.Cells.Item(29, 3).Value = 0.9721315662460447
$ws.Cells.Item(29, 4).Value = 1.097862878032881
$ws.Cells.Item(29, 5).Value = 0.9594025300772112
$ws.Cells.Item(29, 6).Value = 1.097862878032881
$ws.Cells.Item(29, 7).Value = 0.9594025300772112
$ws.Cells.Item(29, 8).Value = 0.8953943628454462
$ws.Cells.Item(29, 9).Value = 1.07810834417008
$ws.Cells.Item(29, 10).Value = 0.9695229094613148
$ws.Cells.Item(29, 11).Value = 0.9594025300772112
$ws.Cells.Item(29, 12).Value = 0.9721315662460447
$ws.Cells.Item(29, 13).Value = 1.034997222139463
$ws.Cells.Item(29, 14).Value = 1.034997222139463
$ws.Cells.Item(29, 15).Value = 1.049367596149668
$ws.Cells.Item(29, 16).Value = 1.009798991452046
$ws.Cells.Item(29, 17).Value = 1.009798991452046
$ws.Cells.Item(29, 18).Value = 0.997199876108337
$ws.Cells.Item(29, 19).Value = 0.997199876108337
$ws.Cells.Item(29, 20).Value = 0.9954037651388296

# Row 30
$ws.Cells.Item(30, 2).Value = 'Michael-CCHex'
$ws.Cells.Item(30, 3).Value = 1.035682238279795
$ws.Cells.Item(30, 4).Value = 0.9243771193540303
$ws.Cells.Item(30, 5).Value = 0.9846569239672783
$ws.Cells.Item(30, 6).Value = 0.9243771193540303
$ws.Cells.Item(30, 7).Value = 0.9846569239672783
$ws.Cells.Item(30, 8).Value = 1.138647962157371
$ws.Cells.Item(30, 9).Value = 0.9432587725932176
$ws.Cells.Item(30, 10).Value = 1.020282325744046
$ws.Cells.Item(30, 11).Value = 0.9846569239672783
$ws.Cells.Item(30, 12).Value = 1.035682238279795
$ws.Cells.Item(30, 13).Value = 0.9800296788169125
$ws.Cells.Item(30, 14).Value = 0.9800296788169125
$ws.Cells.Item(30, 15).Value = 0.9677727100756809
$ws.Cells.Item(30, 16).Value = 0.9815720938670345
$ws.Cells.Item(30, 17).Value = 0.9815720938670344
$ws.Cells.Item(30, 18).Value = 0.9823433013920954
$ws.Cells.Item(30, 19).Value = 0.9823433013920954
$ws.Cells.Item(30, 20).Value = 1.007817557015956

# Row 31
$ws.Cells.Item(31, 2).Value = 'Michael-SNHex'
$ws.Cells.Item(31, 3).Value = 1.165079981003884
$ws.Cells.Item(31, 4).Value = 0.4917599565662167
$ws.Cells.Item(31, 5).Value = 1.031965525543072
$ws.Cells.Item(31, 6).Value = 0.4917599565662167
$ws.Cells.Item(31, 7).Value = 1.031965525543072
$ws.Cells.Item(31, 8).Value = 1.429449040003286
$ws.Cells.Item(31, 9).Value = 0.724330185425365
$ws.Cells.Item(31, 10).Value = 1.125135702788174
$ws.Cells.Item(31, 11).Value = 1.031965525543072
$ws.Cells.Item(31, 12).Value = 1.165079981003884
$ws.Cells.Item(31, 13).Value = 0.8284199687850504
$ws.Cells.Item(31, 14).Value = 0.8284199687850504
$ws.Cells.Item(31, 15).Value = 0.7937233743318219
$ws.Cells.Item(31, 16).Value = 0.896268487704391
$ws.Cells.Item(31, 17).Value = 0.896268487704391
$ws.Cells.Item(31, 18).Value = 0.9301927471640613
$ws.Cells.Item(31, 19).Value = 0.9301927471640613
$ws.Cells.Item(31, 20).Value = 0.9946200652216662

# New trailing rows need the header-style formatting copied onto column A
# (value-only writes above do not carry a cell style).
$ws.Range("A29").Copy()
$ws.Range("A30:A31").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(31, 1).Value = 29